$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.340.18"
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").Value = "1.835.26"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5188"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3234"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06744"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7615"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.73%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.879.44"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07645"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.007"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.48%  "

$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007868"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").Value = "26.389.92"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").Value = "2.067.97"
$ws.Range("E21").Value = "  -3.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.548"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.399"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.924"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.236"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.08%  "

$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.163"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.136"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08701"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04778"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.03%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.856"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6883"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.057"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01757"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.187"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4815"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8877"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.085"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.18%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.643"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4124"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05851"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8833"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
